# Insert two new weekly-report rows into the "Hortaliza, Feria Lagunitas de
# Puerto Montt - Zapallo italiano" sheet at rows 356-357, pushing the
# existing rows 356-371 down to 358-373 (dimension grows from R371 to R373).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 356 (shifts old rows 356:371 -> 358:373)
$ws.Rows.Item(356).Resize(2).Insert()

# New row 356
$ws.Cells.Item(356, 1).Value = 4
$ws.Cells.Item(356, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(356, 3).Value = "Los Lagos"
$ws.Cells.Item(356, 4).Value = 45008
$ws.Cells.Item(356, 5).Value = 10
$ws.Cells.Item(356, 6).Value = 100112032
$ws.Cells.Item(356, 7).Value = "Zapallo italiano"
$ws.Cells.Item(356, 8).Value = "Sin especificar"
$ws.Cells.Item(356, 9).Value = "Primera"
$ws.Cells.Item(356, 10).Value = 70
$ws.Cells.Item(356, 11).Value = 12000
$ws.Cells.Item(356, 12).Value = 12000
$ws.Cells.Item(356, 13).Value = 12000
$ws.Cells.Item(356, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(356, 15).Value = "Región Metropolitana"
$ws.Cells.Item(356, 16).Value = 240
$ws.Cells.Item(356, 17).Value = 50
$ws.Cells.Item(356, 18).Value = "Hortaliza"

# New row 357
$ws.Cells.Item(357, 1).Value = 4
$ws.Cells.Item(357, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(357, 3).Value = "Los Lagos"
$ws.Cells.Item(357, 4).Value = 45008
$ws.Cells.Item(357, 5).Value = 10
$ws.Cells.Item(357, 6).Value = 100112032
$ws.Cells.Item(357, 7).Value = "Zapallo italiano"
$ws.Cells.Item(357, 8).Value = "Sin especificar"
$ws.Cells.Item(357, 9).Value = "Primera"
$ws.Cells.Item(357, 10).Value = 70
$ws.Cells.Item(357, 11).Value = 15000
$ws.Cells.Item(357, 12).Value = 15000
$ws.Cells.Item(357, 13).Value = 15000
$ws.Cells.Item(357, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(357, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(357, 16).Value = 300
$ws.Cells.Item(357, 17).Value = 50
$ws.Cells.Item(357, 18).Value = "Hortaliza"

# Make sure the D356/D357 date cells keep the same number format style (s="2")
# as the rest of the column (the row-insert should already have propagated
# this from the row above, but set explicitly to be safe).
$ws.Cells.Item(356, 4).NumberFormat = $ws.Cells.Item(358, 4).NumberFormat
$ws.Cells.Item(357, 4).NumberFormat = $ws.Cells.Item(358, 4).NumberFormat
